# Sprints.xlsx docu update:
#   - add a "Software development total time" label next to the hours total
#     (new shared string + a header-ish style re-using the existing task
#      description font/fill but without the table border)
#   - leave the selection sitting on that new cell, matching where the
#     author ended up after making the edit
#   - nudge the saved window position, as Excel does when the sheet is
#     rearranged/scrolled before saving

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label cell above/left of the totals formula in D17.
$c17 = $ws.Range("C17")
$c17.Value = "Software development total time"

# Style it like the other "label" cells (dark text on the light fill used
# throughout the sheet), vertically centered, but with no border - this is
# what mints the new cellXfs entry.
$c17.Font.Color = 0
$c17.Interior.Color = 16777215
$c17.VerticalAlignment = -4108

# Leave the selection on the newly added cell.
$c17.Select()

# Saved window position moved slightly to the right.
$excel.Windows.Item(1).Left = 840

Write-Output "Added Software development total time label in C17"
